$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1498.4286
$ws.Range("I29").Value = 97.8
$ws.Range("K29").Value = 293.4
$ws.Range("M29").Value = -12.39999999999998

$ws.Range("H38").Value = 545.8
$ws.Range("J38").Value = 2433.3333
$ws.Range("L38").Value = 7299.999899999999
$ws.Range("N38").Value = -8043.999899999999

$ws.Range("H58").Value = 2065.1428
$ws.Range("I58").Value = 805
$ws.Range("J58").Value = 2408.818
$ws.Range("K58").Value = 2415
$ws.Range("L58").Value = 7226.454000000001
$ws.Range("M58").Value = -2265
$ws.Range("N58").Value = -7526.454000000001

$ws.Range("H64").Value = 3299.932
$ws.Range("I64").Value = 3141.4167
$ws.Range("J64").Value = 3359.375
$ws.Range("K64").Value = 3141.4167
$ws.Range("L64").Value = 3359.375
$ws.Range("M64").Value = -2893.4167
$ws.Range("N64").Value = -3855.375

$ws.Range("H67").Value = 3299.932
$ws.Range("I67").Value = 3141.4167
$ws.Range("J67").Value = 3359.375
$ws.Range("K67").Value = 3141.4167
$ws.Range("L67").Value = 3359.375
$ws.Range("M67").Value = -2283.4167
$ws.Range("N67").Value = -5075.375

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").ClearContents()

$ws.Range("H87").Value = 45313.46
$ws.Range("J87").Value = 45313.46
$ws.Range("L87").Value = 45313.46
$ws.Range("N87").Value = -47809.46

$ws.Range("H90").Value = 45313.46
$ws.Range("J90").Value = 45313.46
$ws.Range("L90").Value = 135940.38
$ws.Range("N90").Value = -148420.38

$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = 4000
$ws.Range("N113").Value = -10508
$ws.Range("L113").ClearContents()

$ws.Range("H135").Value = 699.62964
$ws.Range("I135").Value = 570.0213
$ws.Range("J135").Value = 1569.8572
$ws.Range("K135").Value = 5130.1917
$ws.Range("L135").Value = 14128.7148
$ws.Range("M135").Value = -2595.1917
$ws.Range("N135").Value = -19198.7148

$ws.Range("H138").Value = 2453.6
$ws.Range("I138").Value = 1535.017
$ws.Range("J138").Value = 4538.077
$ws.Range("K138").Value = 4605.051
$ws.Range("L138").Value = 13614.231
$ws.Range("M138").Value = 534.9489999999996
$ws.Range("N138").Value = -23894.231

$ws.Range("H141").Value = 4264.409
$ws.Range("I141").Value = 2028.9744
$ws.Range("J141").Value = 21700.8
$ws.Range("K141").Value = 6086.9232
$ws.Range("L141").Value = 65102.39999999999
$ws.Range("M141").Value = -906.9232000000002
$ws.Range("N141").Value = -75462.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7109.8086
$ws.Range("I32").Value = 7790.037
$ws.Range("J32").Value = 2871.4614
$ws.Range("K32").Value = 7790.037
$ws.Range("L32").Value = 2871.4614
$ws.Range("M32").Value = -7503.037
$ws.Range("N32").Value = -3445.4614

$ws.Range("H61").Value = 1447.4117
$ws.Range("I61").Value = 1289.1875
$ws.Range("J61").Value = 3979
$ws.Range("K61").Value = 1289.1875
$ws.Range("L61").Value = 3979
$ws.Range("M61").Value = -1077.1875
$ws.Range("N61").Value = -4403

$ws.Range("H74").Value = 869.0540999999999
$ws.Range("I74").Value = 813.80646
$ws.Range("J74").Value = 1154.5
$ws.Range("K74").Value = 813.80646
$ws.Range("L74").Value = 1154.5
$ws.Range("M74").Value = 60.19353999999998
$ws.Range("N74").Value = -2902.5

$ws.Range("H77").Value = 869.0540999999999
$ws.Range("I77").Value = 813.80646
$ws.Range("J77").Value = 1154.5
$ws.Range("K77").Value = 4069.0323
$ws.Range("L77").Value = 5772.5
$ws.Range("M77").Value = 298.9677000000001
$ws.Range("N77").Value = -14508.5

$ws.Range("H136").Value = 1447.4117
$ws.Range("I136").Value = 1289.1875
$ws.Range("J136").Value = 3979
$ws.Range("K136").Value = 3867.5625
$ws.Range("L136").Value = 11937
$ws.Range("M136").Value = -1317.5625
$ws.Range("N136").Value = -17037

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2004.7561
$ws.Range("I31").Value = 1413.2333
$ws.Range("J31").Value = 3618
$ws.Range("K31").Value = 1413.2333
$ws.Range("L31").Value = 3618
$ws.Range("M31").Value = -1118.2333
$ws.Range("N31").Value = -4208

$ws.Range("H34").Value = 2004.7561
$ws.Range("I34").Value = 1413.2333
$ws.Range("J34").Value = 3618
$ws.Range("K34").Value = 1413.2333
$ws.Range("L34").Value = 3618
$ws.Range("M34").Value = -1211.2333
$ws.Range("N34").Value = -4022

$ws.Range("H58").Value = 1030243.3
$ws.Range("I58").Value = 1950682.9
$ws.Range("J58").Value = 1516.7059
$ws.Range("K58").Value = 1950682.9
$ws.Range("L58").Value = 1516.7059
$ws.Range("M58").Value = -1950479.9
$ws.Range("N58").Value = -1922.7059

$ws.Range("H62").Value = 57989.332
$ws.Range("I62").Value = 102221
$ws.Range("J62").Value = 2699.75
$ws.Range("K62").Value = 102221
$ws.Range("L62").Value = 2699.75
$ws.Range("M62").Value = -101597
$ws.Range("N62").Value = -3947.75

$ws.Range("H65").Value = 57989.332
$ws.Range("I65").Value = 102221
$ws.Range("J65").Value = 2699.75
$ws.Range("K65").Value = 511105
$ws.Range("L65").Value = 13498.75
$ws.Range("M65").Value = -507985
$ws.Range("N65").Value = -19738.75

$ws.Range("H99").Value = 5900
$ws.Range("I99").Value = 5900
$ws.Range("K99").Value = 5900
$ws.Range("M99").Value = -4402

$ws.Range("H126").Value = 5900
$ws.Range("I126").Value = 5900
$ws.Range("K126").Value = 17700
$ws.Range("M126").Value = -15230

$ws.Range("H132").Value = 339505.88
$ws.Range("I132").Value = 467042.34
$ws.Range("J132").Value = 3273.3635
$ws.Range("K132").Value = 1401127.02
$ws.Range("L132").Value = 9820.0905
$ws.Range("M132").Value = -1398597.02
$ws.Range("N132").Value = -14880.0905

$ws.Range("H134").Value = 2186.074
$ws.Range("I134").Value = 1784.5
$ws.Range("J134").Value = 2770.182
$ws.Range("K134").Value = 5353.5
$ws.Range("L134").Value = 8310.545999999998
$ws.Range("M134").Value = -2818.5
$ws.Range("N134").Value = -13380.546

$ws.Range("H136").Value = 1030243.3
$ws.Range("I136").Value = 1950682.9
$ws.Range("J136").Value = 1516.7059
$ws.Range("K136").Value = 5852048.699999999
$ws.Range("L136").Value = 4550.1177
$ws.Range("M136").Value = -5849498.699999999
$ws.Range("N136").Value = -9650.117699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1194.7097
$ws.Range("I5").Value = 1171.24
$ws.Range("J5").Value = 1292.5
$ws.Range("K5").Value = 3513.72
$ws.Range("L5").Value = 3877.5
$ws.Range("M5").Value = -3401.72
$ws.Range("N5").Value = -4101.5

$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 300
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = -673
$ws.Range("M20").ClearContents()

$ws.Range("H34").Value = 491.64865
$ws.Range("J34").Value = 556.9666999999999
$ws.Range("L34").Value = 1670.9001
$ws.Range("N34").Value = -1838.9001

$ws.Range("H39").Value = 2463.2222
$ws.Range("J39").Value = 2463.2222
$ws.Range("L39").Value = 7389.6666
$ws.Range("N39").Value = -7977.6666

$ws.Range("H55").Value = 5815.3887
$ws.Range("J55").Value = 5815.3887
$ws.Range("L55").Value = 17446.1661
$ws.Range("N55").Value = -17800.1661

$ws.Range("H122").Value = 853.7895
$ws.Range("I122").Value = 500.25
$ws.Range("J122").Value = 1110.909
$ws.Range("K122").Value = 4502.25
$ws.Range("L122").Value = 9998.181
$ws.Range("M122").Value = -2052.25
$ws.Range("N122").Value = -14898.181

$ws.Range("H135").Value = 1194.7097
$ws.Range("I135").Value = 1171.24
$ws.Range("J135").Value = 1292.5
$ws.Range("K135").Value = 10541.16
$ws.Range("L135").Value = 11632.5
$ws.Range("M135").Value = -8006.16
$ws.Range("N135").Value = -16702.5

$ws.Range("H140").Value = 1712.4193
$ws.Range("I140").Value = 1156.2593
$ws.Range("K140").Value = 3468.7779
$ws.Range("M140").Value = 1711.2221

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1541.2858
$ws.Range("I132").Value = 904.2353000000001
$ws.Range("K132").Value = 2712.7059
$ws.Range("M132").Value = -182.7058999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1353.5555
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -505

$ws.Range("H27").Value = 1353.5555
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 800
$ws.Range("M27").Value = -693

$ws.Range("H40").Value = 2904.2
$ws.Range("I40").Value = 2291.7144
$ws.Range("J40").Value = 4333.3335
$ws.Range("K40").Value = 2291.7144
$ws.Range("L40").Value = 4333.3335
$ws.Range("M40").Value = -2155.7144
$ws.Range("N40").Value = -4605.3335

$ws.Range("H122").Value = 11543921
$ws.Range("I122").Value = 9619901
$ws.Range("J122").Value = 15391962
$ws.Range("K122").Value = 28859703
$ws.Range("L122").Value = 46175886
$ws.Range("M122").Value = -28857253
$ws.Range("N122").Value = -46180786

$ws.Range("H132").Value = 6037.625
$ws.Range("I132").Value = 5660.5
$ws.Range("J132").Value = 6666.1665
$ws.Range("K132").Value = 16981.5
$ws.Range("L132").Value = 19998.4995
$ws.Range("M132").Value = -14451.5
$ws.Range("N132").Value = -25058.4995

$ws.Range("H136").Value = 1990.6227
$ws.Range("I136").Value = 1713.3096
$ws.Range("J136").Value = 3049.4546
$ws.Range("K136").Value = 5139.9288
$ws.Range("L136").Value = 9148.363799999999
$ws.Range("M136").Value = -2589.9288
$ws.Range("N136").Value = -14248.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 46000
$ws.Range("J54").Value = 46000
$ws.Range("L54").Value = 46000
$ws.Range("N54").Value = -47040

$ws.Range("H122").Value = 15626422
$ws.Range("I122").Value = 14707166
$ws.Range("K122").Value = 44121498
$ws.Range("M122").Value = -44119048

$ws.Range("H132").Value = 1652.5938
$ws.Range("I132").Value = 1165.0952
$ws.Range("J132").Value = 2583.2727
$ws.Range("K132").Value = 3495.2856
$ws.Range("L132").Value = 7749.8181
$ws.Range("M132").Value = -965.2856000000002
$ws.Range("N132").Value = -12809.8181
